$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.57"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'5.399"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.05760"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Value = "'3.431"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Value = "'6.316"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'0.8105"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "'0.8943"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Value = "'0.1443"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9WazirXWRXBestin24h"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07343"
$ws.Range("D11").Style = "Normal"
$ws.Range("E12").Value = "11BitrueCoinBTR"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.02981"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.09408"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'3.930"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'0.001585"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.04792"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.0005842"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006178"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.004068"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'0.0009946"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Value = "'3.732"
$ws.Range("D23").Style = "Normal"
$ws.Range("D40").Value = "'0.03897"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.006793"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1073"
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").Value = "'0.006791"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005650"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.3801"
$ws.Range("D47").Style = "Normal"
